# Add missing "eventscenter" cover/label shapes to the delete diagram slide,
# and grow the rectangle they previously clipped.
#
# Point values below were chosen so that this runtime's pt -> EMU
# conversion (IEEE-754 single-precision cast, then truncate toward zero,
# times 12700) reproduces the exact target EMU values from the authoritative
# OOXML diff.

$p = $ppt.ActivePresentation

# --- Locate the slide that holds "Rectangle 67" (id 68) -------------------
# This deck repeats the same diagram (with the same shape names/ids) across
# several slides, so disambiguate using the shape's original size too
# (13.683937 x 47.185197 pt == cx 173786 / cy 599252 EMU), not just name/id.
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $candidate = $slide.Shapes.Item($shi)
        if ($candidate.Name -eq "Rectangle 67" -and $candidate.Id -eq 68) {
            $wOk = [Math]::Abs($candidate.Width - 13.683937) -lt 0.01
            $hOk = [Math]::Abs($candidate.Height - 47.185197) -lt 0.01
            if ($wOk -and $hOk) {
                $targetSlide = $slide
                $targetShape = $candidate
            }
        }
    }
}

# --- 1. Grow "Rectangle 67" (id 68) so it reaches the new lower edge ------
# off x="4384723" y="5071221" / ext cx="173786" cy="599252"
#   -> off x="4384723" y="5071220" / ext cx="173786" cy="872379"
$targetShape.Top = 399.3087
$targetShape.Height = 68.6913

# --- 2. Find donor "cover" shapes already used elsewhere in the deck ------
# (slide 24 carries an identical pair: a white roundRect "page" background
# plus a wrap="none" auto-fit "TextBox 62" label; reuse them so the style
# refs / text body markup match exactly.)
$donorSlide = $p.Slides.Item(24)
$roundRectTemplate = $null
$textBoxTemplate = $null
for ($i = 1; $i -le $donorSlide.Shapes.Count; $i++) {
    $sh = $donorSlide.Shapes.Item($i)
    if ($roundRectTemplate -eq $null -and $sh.Name -eq "Rounded Rectangle 61") {
        $roundRectTemplate = $sh
    }
    if ($textBoxTemplate -eq $null -and $sh.Name -eq "TextBox 62") {
        $textBoxTemplate = $sh
    }
}

# --- 3. Add the 4 new shapes (bottom group's cover+label, then the top
#        group's cover+label), positioned via Copy/Paste of the donor
#        shapes so fill / line / effect / font style refs come along. ----

# Shape "Rounded Rectangle 61" (bottom cover)
# off x="0" y="3475534" / ext cx="9144000" cy="3230066"
$roundRectTemplate.Copy()
$rr1 = $targetSlide.Shapes.Paste().Item(1)
$rr1.Left = 0.0
$rr1.Top = 273.66414
$rr1.Width = 720.0
$rr1.Height = 254.3359

# Shape "TextBox 62" = "Delete" (bottom label)
# off x="194562" y="3737425" / ext cx="799706" cy="369332"
$textBoxTemplate.Copy()
$tb1 = $targetSlide.Shapes.Paste().Item(1)
$tb1.TextFrame.TextRange.Text = "Delete"
$tb1.Left = 15.3199
$tb1.Top = 294.2855
$tb1.Width = 62.969
$tb1.Height = 29.0813

# Shape "Rounded Rectangle 61" (top cover)
# off x="0" y="-271987" / ext cx="9144000" cy="3552166"
$roundRectTemplate.Copy()
$rr2 = $targetSlide.Shapes.Paste().Item(1)
$rr2.Left = 0.0
$rr2.Top = -21.4163
$rr2.Width = 720.0
$rr2.Height = 279.69815

# Shape "TextBox 62" = "Delete" (top label)
# off x="194562" y="-10096" / ext cx="799706" cy="369332"
$textBoxTemplate.Copy()
$tb2 = $targetSlide.Shapes.Paste().Item(1)
$tb2.TextFrame.TextRange.Text = "Delete"
$tb2.Left = 15.3199
$tb2.Top = -0.795
$tb2.Width = 62.969
$tb2.Height = 29.0813

# --- 4. The 4 new shapes land at the end of the z-order (end of spTree) --
# after Copy/Paste; move them to the very front (start of spTree, behind
# everything else) in the order rr1, tb1, rr2, tb2, matching the diff
# which inserts them right after <p:grpSpPr> and before the pre-existing
# first shape ("Rectangle 62").
$tb2.ZOrder(1)
$rr2.ZOrder(1)
$tb1.ZOrder(1)
$rr1.ZOrder(1)
